$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.385465741157532
$ws.Range("B1").Value = 2.664197683334351
$ws.Range("C1").Value = 5.914616584777832
$ws.Range("D1").Value = 2.281152009963989
$ws.Range("E1").Value = 1.207550287246704
